# Daily attendance processing - 2026-01-08 08:43:19
#
# Normalizes the "Recorded By" (column G) cell text on the active sheet.
# Several rows list the same recorders in different orders (e.g. the
# automated "System" user sometimes gets logged before, sometimes after
# the human recorder). This pass re-orders a known set of "Recorded By"
# combinations to a canonical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Canonical re-ordering map: old combined text -> new combined text.
$map = @{
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value2

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
